$d = $word.ActiveDocument

# 1. Locate the "Assignation d'utilisateur..." list paragraph (currently ends the list
#    with a period) and change its final period to a comma, since two more list items
#    will be appended after it.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Assignation d")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Assignation d'utilisateur...' list paragraph."
}

$target = $d.Paragraphs($targetIndex)
$periodRange = $d.Range($target.Range.End - 2, $target.Range.End - 1)
$periodRange.Text = ","

# 2. Insert two new list paragraphs after it, inheriting the same list/paragraph
#    formatting (Paragraphedeliste style, numId 3 bullet list, Verdana font).
$target.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs($targetIndex + 1)
$newPara1.Range.Text = "Création et gestion de tags pour les tâches,"

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs($targetIndex + 2)
$newPara2.Range.Text = "Commenter sur les tâches."
